$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New expense rows (row 32-50) for "Nomina" payroll entries from Feb 21/27 2026
# and misc categories (Insumos, Servicios, Otros) loaded in bulk.
$rows = @(
  ,@(1771743882115, "2026-02-21", "Martha", "Nomina", 80000)
  ,@(1771743914552, "2026-02-21", "Dreisy", "Nomina", 80000)
  ,@(1771743947518, "2026-02-21", "Luisa vallejo", "Nomina", 80000)
  ,@(1771743955114, "2026-02-21", "Julian", "Nomina", 80000)
  ,@(1771743971143, "2026-02-21", "Alexander", "Nomina", 70000)
  ,@(1771744008831, "2026-02-21", "Juanita", "Nomina", 70000)
  ,@(1771744022489, "2026-02-21", "Valeria", "Nomina", 70000)
  ,@(1771744045848, "2026-02-21", "Jhojan", "Nomina", 70000)
  ,@(1771744059511, "2026-02-21", "Sebastian", "Nomina", 150000)
  ,@(1771744073710, "2026-02-21", "Juan Carlos", "Nomina", 80000)
  ,@(1771744099156, "2026-02-21", "Laura", "Nomina", 70000)
  ,@(1771744302410, "2026-02-21", "Papel, trapos, cloro", "Insumos", 42370)
  ,@(1772230354842, "2026-02-27", "Energia", "Servicios", 723540)
  ,@(1772230494337, "2026-02-27", "Prestamo Mr Tango- Sayci-Acinpro", "Otros", 5800000)
  ,@(1772234485905, "2026-02-27", "Zumo de limon", "Insumos", 40000)
  ,@(1772246917060, "2026-02-27", "Jhojan Buitrago", "Nomina", 80000)
  ,@(1772246944016, "2026-02-27", "Fredy Ramirez", "Nomina", 80000)
  ,@(1772247070113, "2026-02-27", "Sebastian ", "Nomina", 80000)
  ,@(1772247096948, "2026-02-27", "Luisa Vallejo", "Nomina", 80000)
)

$responsable = '{"nombre":"Luisa","rol":"ADMIN"}'

$startRow = 32
for ($i = 0; $i -lt $rows.Count; $i++) {
  $r = $startRow + $i
  $data = $rows[$i]

  # B:D and F hold text values in the source data (ids/amounts stay numeric).
  # Force text NumberFormat first so Excel does not auto-coerce date-like
  # strings (e.g. "2026-02-21") into date serials, then clear the format
  # back off so the cell keeps the default (General) style, matching the
  # rest of the sheet.
  $ws.Range("B$r`:D$r").NumberFormat = "@"
  $ws.Range("F$r").NumberFormat = "@"

  $ws.Cells.Item($r, 1).Value = $data[0]
  $ws.Cells.Item($r, 2).Value = $data[1]
  $ws.Cells.Item($r, 3).Value = $data[2]
  $ws.Cells.Item($r, 4).Value = $data[3]
  $ws.Cells.Item($r, 5).Value = $data[4]
  $ws.Cells.Item($r, 6).Value = $responsable

  $ws.Range("A$r`:F$r").ClearFormats()
}
